$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pairing text for the "winners" column (B), rows 2-12 introduce the 11
# distinct team names (in this exact order so the shared-string table is
# rebuilt with matching indices); rows 13-17 repeat some of those teams.
$teams = @(
    "הקבוצה של: אברג'ל, דוידזון, מודן",
    "הקבוצה של: גל ארצי, בן שושן",
    "הקבוצה של: אריאל בן אליעזר, ליאורה",
    "הקבוצה של: תמיר, הרצברג",
    "הקבוצה של: גל נימצקי, עמרי נצן",
    "הקבוצה של: עמית לוי, שקד",
    "הקבוצה של: פייביש, דור פרידמן",
    "הקבוצה של: יואב סטרולוביץ', עמית בר",
    "הקבוצה של: ירין, מנחה",
    "הקבוצה של: עמרי קונסטנטינו, איילה",
    "הקבוצה של: אריאל ליבזון, הוד"
)

$rowTeamIndex = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 3
    6  = 4
    7  = 5
    8  = 6
    9  = 7
    10 = 8
    11 = 9
    12 = 10
    13 = 0
    14 = 2
    15 = 5
    16 = 4
    17 = 7
}

foreach ($row in 2..17) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $teams[$rowTeamIndex[$row]]
    $cell.HorizontalAlignment = -4152
    $cell.WrapText = $true
}

# Rows 18-27 no longer hold a pairing result; clear their content (keeping
# the existing cell formatting).
foreach ($row in 18..27) {
    $ws.Cells.Item($row, 2).ClearContents()
}

# The active view now shows row 13 at the top of the frozen pane, with the
# selection on D17.
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("A13").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A13").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D17").Select()
